$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.790.29"
$ws.Range("E2").Value = "  -1.86%  "

# Row 3
$ws.Range("D3").Value = "'1.545.55"
$ws.Range("E3").Value = "  -1.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'204.44"
$ws.Range("E5").Value = "  -1.74%  "

# Row 6
$ws.Range("D6").Value = "'0.481"
$ws.Range("E6").Value = "  -1.83%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.245"
$ws.Range("E8").Value = "  -1.27%  "

# Row 9
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'21.33"
$ws.Range("E9").Value = "  -4.51%  "

# Row 10
$ws.Range("E10").Value = "  -1.91%  "

# Row 11
$ws.Range("E11").Value = "  -1.22%  "

# Row 12
$ws.Range("D12").Value = "'1.765.89"
$ws.Range("E12").Value = "  -1.86%  "

# Row 13
$ws.Range("D13").Value = "'1.544.96"
$ws.Range("E13").Value = "  -2.07%  "

# Row 14
$ws.Range("E14").Value = "  -2.90%  "

# Row 15
$ws.Range("E15").Value = "  -2.32%  "

# Row 16
$ws.Range("D16").Value = "'26.779.47"
$ws.Range("E16").Value = "  -1.86%  "

# Row 17
$ws.Range("D17").Value = "'60.91"
$ws.Range("E17").Value = "  -2.83%  "

# Row 18
$ws.Range("D18").Value = "'213.63"
$ws.Range("E18").Value = "  -1.21%  "

# Row 19
$ws.Range("D19").Value = "'7.24"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0680"
$ws.Range("E20").Value = "  -1.19%  "

# Row 21
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("E22").Value = "  -1.72%  "

# Row 23
$ws.Range("D23").Value = "'9.02"
$ws.Range("E23").Value = "  -4.12%  "

# Row 24
$ws.Range("D24").Value = "'2.00"
$ws.Range("E24").Value = "  -0.36%  "

# Row 25
$ws.Range("D25").Value = "'152.32"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26
$ws.Range("D26").Value = "'6.49"
$ws.Range("E26").Value = "  -2.95%  "

# Row 27
$ws.Range("D27").Value = "'14.79"
$ws.Range("E27").Value = "  -1.18%  "

# Row 28
$ws.Range("E28").Value = "  +0.18%  "

# Row 29
$ws.Range("E29").Value = "  -2.55%  "

# Row 30
$ws.Range("D30").Value = "'0.0460"
$ws.Range("E30").Value = "  -0.83%  "

# Row 31
$ws.Range("E31").Value = "  -3.36%  "

# Row 32
$ws.Range("E32").Value = "  -0.76%  "

# Row 33
$ws.Range("D33").Value = "'1.357.77"
$ws.Range("E33").Value = "  -3.58%  "

# Row 34
$ws.Range("E34").Value = "  -1.16%  "

# Row 35
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  -4.99%  "

# Row 36
$ws.Range("E36").Value = "  -0.63%  "

# Row 37
$ws.Range("D37").Value = "'0.911"
$ws.Range("E37").Value = "  -3.34%  "

# Row 38
$ws.Range("E38").Value = "  -2.37%  "

# Row 39
$ws.Range("E39").Value = "  +0.61%  "

# Row 40
$ws.Range("E40").Value = "  -2.58%  "

# Row 41
$ws.Range("E41").Value = "  +0.16%  "

# Row 42
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  -0.85%  "

# Row 43
$ws.Range("E43").Value = "  +2.64%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").Value = "'1.76"
$ws.Range("E45").Value = "  -2.69%  "

# Row 46
$ws.Range("D46").Value = "'62.75"
$ws.Range("E46").Value = "  -1.92%  "

# Row 47
$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  -2.96%  "

# Row 48
$ws.Range("D48").Value = "'1.679.75"
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("D49").Value = "'85.73"
$ws.Range("E49").Value = "  -0.71%  "

# Row 50
$ws.Range("D50").Value = "'0.0506"
$ws.Range("E50").Value = "  +2.46%  "

# Row 51
$ws.Range("D51").Value = "'0.0₇0962"
$ws.Range("E51").Value = "  -2.71%  "

